# Update NATMI ligand-receptor TPM-derived values (Mdk-Lrp1) per new TPM recomputation.
# Ligand/receptor average & total expression (and all downstream specificity /
# edge-weight columns that are derived from them) were refreshed for the new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "G2" = 0.5979736666666666
    "H2" = 1.793921
    "I2" = 0.03342655292740804
    "J2" = 0.03342655292740804
    "M2" = 3.456265333333333
    "N2" = 10.368796
    "O2" = 0.009841535807677501
    "P2" = 0.0098415358076775
    "Q2" = 2.066755654346222
    "R2" = 18.600800889116
    "S2" = 0.0003289686175623135
    "T2" = 0.0003289686175623134
    "G3" = 0.5979736666666666
    "H3" = 1.793921
    "I3" = 0.03342655292740804
    "J3" = 0.03342655292740804
    "O3" = 0.8587907398420774
    "P3" = 0.8587907398420773
    "Q3" = 180.3489467654184
    "R3" = 1623.140520888766
    "S3" = 0.02870641411889911
    "T3" = 0.02870641411889911
    "G4" = 0.5979736666666666
    "H4" = 1.793921
    "I4" = 0.03342655292740804
    "J4" = 0.03342655292740804
    "O4" = 0.1313677243502452
    "P4" = 0.1313677243502452
    "Q4" = 27.58766440575878
    "R4" = 248.288979651829
    "S4" = 0.004391170190946621
    "T4" = 0.00439117019094662
    "I5" = 0.8874158839838097
    "J5" = 0.8874158839838097
    "M5" = 3.456265333333333
    "N5" = 10.368796
    "O5" = 0.009841535807677501
    "P5" = 0.0098415358076775
    "Q5" = 54.86870871678622
    "R5" = 493.818378451076
    "S5" = 0.008733535198528446
    "T5" = 0.008733535198528444
    "I6" = 0.8874158839838097
    "J6" = 0.8874158839838097
    "O6" = 0.8587907398420774
    "P6" = 0.8587907398420773
    "S6" = 0.7621045435540671
    "T6" = 0.762104543554067
    "I7" = 0.8874158839838097
    "J7" = 0.8874158839838097
    "O7" = 0.1313677243502452
    "P7" = 0.1313677243502452
    "S7" = 0.1165778052312143
    "T7" = 0.1165778052312142
    "I8" = 0.07915756308878232
    "J8" = 0.07915756308878232
    "M8" = 3.456265333333333
    "N8" = 10.368796
    "O8" = 0.009841535807677501
    "P8" = 0.0098415358076775
    "Q8" = 4.894292912981333
    "R8" = 44.048636216832
    "S8" = 0.0007790319915867421
    "T8" = 0.0007790319915867419
    "I9" = 0.07915756308878232
    "J9" = 0.07915756308878232
    "O9" = 0.8587907398420774
    "P9" = 0.8587907398420773
    "Q9" = 427.0851129215147
    "R9" = 3843.766016293632
    "S9" = 0.06797978216911128
    "T9" = 0.06797978216911127
    "I10" = 0.07915756308878232
    "J10" = 0.07915756308878232
    "O10" = 0.1313677243502452
    "P10" = 0.1313677243502452
    "Q10" = 65.33046618397869
    "R10" = 587.9741956558081
    "S10" = 0.0103987489280843
    "T10" = 0.0103987489280843
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
